$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.457.00'
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").Value = '2.899.27'
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.75'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.70'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.98%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.500'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.42%  '
$ws.Range("D9").Value = '2.899.61'
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.90'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.149'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.431'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.89%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000236'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.71'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("D16").Value = '3.384.05'
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").Value = '62.494.25'
$ws.Range("E17").Value = '  +1.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.59'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").Value = '2.898.08'
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '425.61'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.03'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.657'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.85'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.62'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.86'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.03'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.14%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.02'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000109'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.06'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.48'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.01'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.26%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.71'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.105'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.951'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.36'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.92'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '48.52'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.89'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -5.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '41.25'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +5.83%  '
$ws.Range("E42").Value = '  -1.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.03'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.266'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.51%  '
$ws.Range("D45").Value = '2.712.46'
$ws.Range("E45").Value = '  +0.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0337'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '132.58'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '355.45'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.75%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000217'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +12.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.102'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.06%  '
